# Regenerate the linear/quadratic problem data (new random coefficients) on
# the "Restricciones_del_follower", "Punto_modificado", "Vector_bf" and
# "Vector_BF" sheets, matching the target commit
# "volver a generar problemas cuadraticos y lineales".
#
# Sheets are addressed by index because "Vector_bf" / "Vector_BF" only
# differ by case and Worksheets.Item(name) resolves case-insensitively.

$wb = $excel.ActiveWorkbook

# Helper: write a value that must be stored as TEXT even though it looks
# like a plain number (matches the source file, where every one of these
# cells is a shared string, not a numeric cell). Force text format, write
# the value, then strip the formatting back off so no stray cell style is
# left applied to the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# --- Sheet 3: Restricciones_del_follower ---------------------------------
$wsFollower = $wb.Worksheets.Item(3)

$wsFollower.Range("A2").Value = "6.915 - x - 0.5y"
Set-TextValue $wsFollower.Range("B2") "-4.915"
Set-TextValue $wsFollower.Range("D2") "0.21"
Set-TextValue $wsFollower.Range("E2") "9.9"
Set-TextValue $wsFollower.Range("F2") "5.300000000000001"

$wsFollower.Range("A3").Value = "-1.0500000000000003 - 0.25x + y"
Set-TextValue $wsFollower.Range("B3") "-0.9499999999999997"
Set-TextValue $wsFollower.Range("D3") "0.19"
Set-TextValue $wsFollower.Range("E3") "8.9"
Set-TextValue $wsFollower.Range("F3") "8.9"

$wsFollower.Range("A4").Value = "-6.915 + x + 0.5y"
Set-TextValue $wsFollower.Range("B4") "-1.085"
Set-TextValue $wsFollower.Range("D4") "0.97"
Set-TextValue $wsFollower.Range("E4") "1.0"
Set-TextValue $wsFollower.Range("F4") "1.7000000000000002"

$wsFollower.Range("A5").Value = "-3.2600000000000007 + x - 2y"
Set-TextValue $wsFollower.Range("B5") "-1.2600000000000007"
Set-TextValue $wsFollower.Range("D5") "0.43"
Set-TextValue $wsFollower.Range("E5") "4.6000000000000005"
Set-TextValue $wsFollower.Range("F5") "4.3"

$wsFollower.Range("A6").Value = "-2.47 - y"
Set-TextValue $wsFollower.Range("B6") "-2.47"
Set-TextValue $wsFollower.Range("D6") "0.47"
Set-TextValue $wsFollower.Range("E6") "2.5"
Set-TextValue $wsFollower.Range("F6") "8.4"

# --- Sheet 4: Punto_modificado --------------------------------------------
$wsPunto = $wb.Worksheets.Item(4)
Set-TextValue $wsPunto.Range("A2") "5.68"
Set-TextValue $wsPunto.Range("B2") "2.47"

# --- Sheet 5: Vector_bf ----------------------------------------------------
$wsVecbf = $wb.Worksheets.Item(5)
Set-TextValue $wsVecbf.Range("A2") "1.76"

# --- Sheet 6: Vector_BF -----------------------------------------------------
$wsVecBF = $wb.Worksheets.Item(6)
Set-TextValue $wsVecBF.Range("A2") "5.5249999999999995"
Set-TextValue $wsVecBF.Range("A3") "6.250000000000001"
